# Insert a new data row at row 11 (pushes existing rows 11..52 down to 12..53)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly price observation
$ws.Cells.Item(11, 1).Value = 11
$ws.Cells.Item(11, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(11, 3).Value = "Bíobío"
$ws.Cells.Item(11, 4).Value = 45069
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 100114007
$ws.Cells.Item(11, 7).Value = "Jengibre"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 30
$ws.Cells.Item(11, 11).Value = 15000
$ws.Cells.Item(11, 12).Value = 15000
$ws.Cells.Item(11, 13).Value = 15000
$ws.Cells.Item(11, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(11, 15).Value = "Perú"
$ws.Cells.Item(11, 16).Value = 1154
$ws.Cells.Item(11, 17).Value = 13
$ws.Cells.Item(11, 18).Value = "Hortaliza"
